$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 63, shifting existing rows 63-162 down to 65-164.
# This also copies formatting (e.g. the date style on column D) from the row above,
# matching the original file's style (s="2") on column D cells.
$ws.Rows.Item(63).Resize(2).Insert()

# Fill in the two new rows (63 and 64) with the new daily records.
# Columns A,B,C,E,F,G,H,I,R are constant across all data rows in this sheet.

# Row 63
$ws.Cells.Item(63,1).Value  = 10
$ws.Cells.Item(63,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(63,3).Value  = "La Araucanía"
$ws.Cells.Item(63,4).Value  = 44413
$ws.Cells.Item(63,5).Value  = 9
$ws.Cells.Item(63,6).Value  = 100112037
$ws.Cells.Item(63,7).Value  = "Cebollín"
$ws.Cells.Item(63,8).Value  = "Sin especificar"
$ws.Cells.Item(63,9).Value  = "Primera"
$ws.Cells.Item(63,10).Value = 70
$ws.Cells.Item(63,11).Value = 8000
$ws.Cells.Item(63,12).Value = 8000
$ws.Cells.Item(63,13).Value = 8000
$ws.Cells.Item(63,14).Value = "$/docena de paquetes"
$ws.Cells.Item(63,15).Value = "Provincia de Cautín"
$ws.Cells.Item(63,16).Value = 667
$ws.Cells.Item(63,17).Value = 12
$ws.Cells.Item(63,18).Value = "Hortaliza"

# Row 64
$ws.Cells.Item(64,1).Value  = 10
$ws.Cells.Item(64,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(64,3).Value  = "La Araucanía"
$ws.Cells.Item(64,4).Value  = 44413
$ws.Cells.Item(64,5).Value  = 9
$ws.Cells.Item(64,6).Value  = 100112037
$ws.Cells.Item(64,7).Value  = "Cebollín"
$ws.Cells.Item(64,8).Value  = "Sin especificar"
$ws.Cells.Item(64,9).Value  = "Primera"
$ws.Cells.Item(64,10).Value = 30
$ws.Cells.Item(64,11).Value = 5000
$ws.Cells.Item(64,12).Value = 5000
$ws.Cells.Item(64,13).Value = 5000
$ws.Cells.Item(64,14).Value = "$/docena de paquetes"
$ws.Cells.Item(64,15).Value = "Región de O'Higgins"
$ws.Cells.Item(64,16).Value = 417
$ws.Cells.Item(64,17).Value = 12
$ws.Cells.Item(64,18).Value = "Hortaliza"
